# Applies the crypto price/volume update described by the commit diff.
# Column D holds numeric-looking text (e.g. "23.657.56", "0.3770") that must
# stay plain text (matching the original inlineStr cells) rather than being
# auto-coerced to a number by Excel's input parser (which would also strip
# trailing zeros / reformat the value). We write those with a leading
# apostrophe to force text entry, then restore the cell's original Style so
# no stray "quote prefix" formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range("D2").Style()
$ws.Range("D2").Value = "'23.657.56"
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = "  -1.69%  "
$style = $ws.Range("D3").Style()
$ws.Range("D3").Value = "'1.625.48"
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = "  -2.11%  "
$ws.Range("E4").Value = "  +0.49%  "
$ws.Range("B5").Value = "USDC"
$ws.Range("C5").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$style = $ws.Range("D5").Style()
$ws.Range("D5").Value = "'1.003"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$style = $ws.Range("D6").Style()
$ws.Range("D6").Value = "'307.07"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  -1.00%  "
$style = $ws.Range("D7").Style()
$ws.Range("D7").Value = "'0.3826"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  -1.93%  "
$style = $ws.Range("D8").Style()
$ws.Range("D8").Value = "'0.3770"
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = "  -2.79%  "
$style = $ws.Range("D9").Style()
$ws.Range("D9").Value = "'50.04"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = "  -2.72%  "
$style = $ws.Range("D10").Style()
$ws.Range("D10").Value = "'1.309"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  -4.46%  "
$style = $ws.Range("D11").Style()
$ws.Range("D11").Value = "'1.004"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  +0.45%  "
$style = $ws.Range("D12").Style()
$ws.Range("D12").Value = "'0.08303"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  -2.55%  "
$style = $ws.Range("D13").Style()
$ws.Range("D13").Value = "'23.48"
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = "  -2.20%  "
$style = $ws.Range("D14").Style()
$ws.Range("D14").Value = "'6.791"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  -6.03%  "
$style = $ws.Range("D15").Style()
$ws.Range("D15").Value = "'7.645"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  -5.27%  "
$style = $ws.Range("D16").Style()
$ws.Range("D16").Value = "'0.00001281"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  -2.62%  "
$style = $ws.Range("D17").Style()
$ws.Range("D17").Value = "'1.642.71"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  -0.96%  "
$style = $ws.Range("D18").Style()
$ws.Range("D18").Value = "'92.78"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  -1.96%  "
$style = $ws.Range("D19").Style()
$ws.Range("D19").Value = "'0.06922"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  -0.94%  "
$style = $ws.Range("D20").Style()
$ws.Range("D20").Value = "'19.05"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  -4.75%  "
$style = $ws.Range("D21").Style()
$ws.Range("D21").Value = "'6.788"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  -2.90%  "
$ws.Range("E22").Value = "  +0.19%  "
$style = $ws.Range("D23").Style()
$ws.Range("D23").Value = "'13.40"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  -2.22%  "
$style = $ws.Range("D24").Style()
$ws.Range("D24").Value = "'23.668.05"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  -1.62%  "
$style = $ws.Range("D25").Style()
$ws.Range("D25").Value = "'2.419"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  -2.95%  "
$style = $ws.Range("D26").Style()
$ws.Range("D26").Value = "'2.837"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  -8.66%  "
$style = $ws.Range("D27").Style()
$ws.Range("D27").Value = "'21.68"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  -2.70%  "
$style = $ws.Range("D28").Style()
$ws.Range("D28").Value = "'151.88"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  -1.38%  "
$style = $ws.Range("D29").Style()
$ws.Range("D29").Value = "'5.423"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  +2.00%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$style = $ws.Range("D30").Style()
$ws.Range("D30").Value = "'7.933"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  -0.73%  "
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$style = $ws.Range("D31").Style()
$ws.Range("D31").Value = "'134.85"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  -4.09%  "
$style = $ws.Range("D32").Style()
$ws.Range("D32").Value = "'2.481"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  -0.52%  "
$style = $ws.Range("D33").Style()
$ws.Range("D33").Value = "'1.800.98"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  -2.05%  "
$style = $ws.Range("D34").Style()
$ws.Range("D34").Value = "'0.9733"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  -7.44%  "
$style = $ws.Range("D35").Style()
$ws.Range("D35").Value = "'0.07741"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  -5.09%  "
$style = $ws.Range("D36").Style()
$ws.Range("D36").Value = "'0.02850"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  -5.30%  "
$style = $ws.Range("D37").Style()
$ws.Range("D37").Value = "'6.500"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  -3.70%  "
$style = $ws.Range("D38").Style()
$ws.Range("D38").Value = "'0.2614"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  -3.57%  "
$style = $ws.Range("D39").Style()
$ws.Range("D39").Value = "'10.26"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  -8.43%  "
$style = $ws.Range("D40").Style()
$ws.Range("D40").Value = "'0.08985"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  -1.95%  "
$style = $ws.Range("D41").Style()
$ws.Range("D41").Value = "'0.7370"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  -2.95%  "
$style = $ws.Range("D42").Style()
$ws.Range("D42").Value = "'13.17"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  -3.85%  "
$style = $ws.Range("D43").Style()
$ws.Range("D43").Value = "'1.402"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  -1.48%  "
$style = $ws.Range("D44").Style()
$ws.Range("D44").Value = "'16.33"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  -1.57%  "
$style = $ws.Range("D45").Style()
$ws.Range("D45").Value = "'0.6788"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  -3.58%  "
$style = $ws.Range("D46").Style()
$ws.Range("D46").Value = "'2.380"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  -5.08%  "
$style = $ws.Range("D47").Style()
$ws.Range("D47").Value = "'4.050"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  -1.18%  "
$ws.Range("E48").Value = "  +0.40%  "
$style = $ws.Range("D49").Style()
$ws.Range("D49").Value = "'0.08141"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  -1.97%  "
$style = $ws.Range("D50").Style()
$ws.Range("D50").Value = "'132.64"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  -2.36%  "
$style = $ws.Range("D51").Style()
$ws.Range("D51").Value = "'1.202"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  -3.33%  "
